$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (rows whose computed height shifted slightly,
#     e.g. thick-bottom-border rows and the wrapped-text header row) ---
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 84.5
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 15
$ws.Rows.Item(37).RowHeight = 15
$ws.Rows.Item(38).RowHeight = 15
$ws.Rows.Item(44).RowHeight = 15
$ws.Rows.Item(45).RowHeight = 15
$ws.Rows.Item(60).RowHeight = 15
$ws.Rows.Item(61).RowHeight = 15
$ws.Rows.Item(62).RowHeight = 15
$ws.Rows.Item(63).RowHeight = 15
$ws.Rows.Item(66).RowHeight = 15
$ws.Rows.Item(75).RowHeight = 15
$ws.Rows.Item(79).RowHeight = 15
$ws.Rows.Item(83).RowHeight = 15
$ws.Rows.Item(90).RowHeight = 15
$ws.Rows.Item(96).RowHeight = 15
$ws.Rows.Item(98).RowHeight = 15
$ws.Rows.Item(106).RowHeight = 15.5

# --- Fix mojibake text in shared string used by footnote cell A103 ---
# (Regional Economic Communities footnote: "PALOP" and "MERCOSUR" definitions
#  had UTF-8 bytes that had been mis-decoded as Latin-1/cp1252; restore correct text.)
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'

# --- Update recomputed aggregate/group statistics (columns C-G) for several summary rows ---
# Row 63
$ws.Range("C63").Value = 2.59
$ws.Range("D63").Value = 6.3070000000000004
$ws.Range("E63").Value = 21.552
$ws.Range("F63").Value = 34.343000000000004
$ws.Range("G63").Value = 41.953000000000003

# Row 64
$ws.Range("C64").Value = 3.4772729999999998
$ws.Range("D64").Value = 6.8363639999999997
$ws.Range("E64").Value = 23.290908999999999
$ws.Range("F64").Value = 43.3
$ws.Range("G64").Value = 49.095455000000001

# Row 65
$ws.Range("C65").Value = 5.2272730000000003
$ws.Range("D65").Value = 15.022727
$ws.Range("E65").Value = 48.159090999999997
$ws.Range("F65").Value = 32.831817999999998
$ws.Range("G65").Value = 41.25

# Row 66
$ws.Range("C66").Value = 13.066667000000001
$ws.Range("D66").Value = 20.331250000000001
$ws.Range("E66").Value = 38.925694
$ws.Range("F66").Value = 36.115971999999999
$ws.Range("G66").Value = 43.498610999999997

# Row 76
$ws.Range("C76").Value = 5.4714289999999997
$ws.Range("D76").Value = 15.057143
$ws.Range("E76").Value = 45.228571000000002
$ws.Range("F76").Value = 36.285713999999999
$ws.Range("G76").Value = 44.057143000000003

# Row 79
$ws.Range("C79").Value = 0.72432399999999997
$ws.Range("D79").Value = 1.208108
$ws.Range("E79").Value = 3.9432429999999998
$ws.Range("F79").Value = 33.294595000000001
$ws.Range("G79").Value = 40.948649000000003

# Row 83
$ws.Range("C83").Value = 2.7155559999999999
$ws.Range("D83").Value = 6.621111
$ws.Range("E83").Value = 21.915555999999999
$ws.Range("F83").Value = 34.367778000000001
$ws.Range("G83").Value = 41.98

# Row 89
$ws.Range("C89").Value = 2.221622
$ws.Range("D89").Value = 5.4486489999999996
$ws.Range("E89").Value = 25.389188999999998
$ws.Range("F89").Value = 36.678378000000002
$ws.Range("G89").Value = 43.981081000000003

# Row 90
$ws.Range("C90").Value = 0.56304299999999996
$ws.Range("D90").Value = 0.91521699999999995
$ws.Range("E90").Value = 3.0586959999999999
$ws.Range("F90").Value = 32.808695999999998
$ws.Range("G90").Value = 40.432608999999999
